# "Merge tra versioni Stefano e David"
# Reconcile the two authors' edits of the Misuratori/Punti di misura sheet:
#  - refresh a handful of GPS coordinates (lat/lon) to rounder, re-surveyed values
#  - rename / re-label several measurement points
#  - change the status of one point to "Inacessibile"
#  - add a brand-new row (14) for "I salto Acquedotto Merone" with
#    sexagesimal-to-decimal coordinate formulas
#  - tidy up the sheet view (zoom, selection) and the trailing spacer row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Partitore Trebisaccce - re-surveyed coordinates -------------
$ws.Range("F2").Value = 39.860277777777782
$ws.Range("G2").Value = 16.508611111111112

# --- Row 3: renamed from "I salto Acquedotto Merone (camera di manovra)"
#     to "Camera di manovra Merone", re-surveyed coordinates -------------
$ws.Range("A3").Value = "Camera di manovra Merone"
$ws.Range("F3").Value = 39.176666666666669
$ws.Range("G3").Value = 16.341666666666665

# --- Row 10: Acquedotto Zumpo - re-surveyed coordinates -----------------
$ws.Range("F10").Value = 39.240833333333335
$ws.Range("G10").Value = 16.361944444444447

# --- Row 11: renamed from "SA" to "Sorgente Zumpo", re-surveyed coords --
$ws.Range("A11").Value = "Sorgente Zumpo"
$ws.Range("F11").Value = 39.231111111111112
$ws.Range("G11").Value = 16.404999999999998

# --- Row 12: renamed from "Sorical Differenzi Partitore Musco 2" to "SA",
#     re-surveyed coordinates --------------------------------------------
$ws.Range("A12").Value = "SA"
$ws.Range("F12").Value = 38.999938888888892
$ws.Range("G12").Value = 17.062283333333333

# --- Row 13: renamed from "Sorical Pisarello" to
#     "Sorical Differenzi Murate / Pisarello", re-surveyed coordinates ---
$ws.Range("A13").Value = "Sorical Differenzi Murate / Pisarello"
$ws.Range("F13").Value = 39.116572222222224
$ws.Range("G13").Value = 16.749461111111113

# --- Row 9: status changed to "Inacessibile" ----------------------------
$ws.Range("K9").Value = "Inacessibile"

# --- Remove the old bottom spacer row (was r=23) and add the new data
#     row 14, then re-add the spacer one row higher (r=22) --------------
$ws.Rows("23").Delete()

$ws.Range("A14").Value = "I salto Acquedotto Merone"
$ws.Range("F14").Formula = "=39+10/60+55/3600"
$ws.Range("G14").Formula = "=16+20/60+9/3600"
$ws.Range("J14").Value = "areatecnica"
$ws.Range("K14").Value = "In valutazione"

$ws.Rows("22").RowHeight = 9

# --- Column K (Stato) widened to fit content ----------------------------
$ws.Columns("K").ColumnWidth = 28.8

# --- Sheet view: zoom to 145%, scroll so column C is leftmost, select L10
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
$ws.Range("L10").Select() | Out-Null
